# Update Ccl28-Ccr3 LR-pair sheet with refreshed TPM-based NATMI values.
# The "ECs" cluster rows are replaced by "FAPs" -> Resolving-Mac and
# "Resolving-Mac" -> Resolving-Mac rows, and the old 4-row table collapses
# to 2 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs | Ccl28 | Ccr3 | Resolving-Mac | ...
$ws.Cells.Item(2,1).Value  = "FAPs"
$ws.Cells.Item(2,2).Value  = "Ccl28"
$ws.Cells.Item(2,3).Value  = "Ccr3"
$ws.Cells.Item(2,4).Value  = "Resolving-Mac"
$ws.Cells.Item(2,5).Value  = 1
$ws.Cells.Item(2,6).Value  = 0.3333333333333333
$ws.Cells.Item(2,7).Value  = 0.07943533333333333
$ws.Cells.Item(2,8).Value  = 0.238306
$ws.Cells.Item(2,9).Value  = 0.8002511845635669
$ws.Cells.Item(2,10).Value = 0.8002511845635669
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.1790523333333333
$ws.Cells.Item(2,14).Value = 0.537157
$ws.Cells.Item(2,15).Value = 1
$ws.Cells.Item(2,16).Value = 1
$ws.Cells.Item(2,17).Value = 0.01422308178244444
$ws.Cells.Item(2,18).Value = 0.128007736042
$ws.Cells.Item(2,19).Value = 0.8002511845635669
$ws.Cells.Item(2,20).Value = 0.8002511845635669

# Row 3: Resolving-Mac | Ccl28 | Ccr3 | Resolving-Mac | ...
$ws.Cells.Item(3,1).Value  = "Resolving-Mac"
$ws.Cells.Item(3,2).Value  = "Ccl28"
$ws.Cells.Item(3,3).Value  = "Ccr3"
$ws.Cells.Item(3,4).Value  = "Resolving-Mac"
$ws.Cells.Item(3,5).Value  = 1
$ws.Cells.Item(3,6).Value  = 0.3333333333333333
$ws.Cells.Item(3,7).Value  = 0.01982766666666667
$ws.Cells.Item(3,8).Value  = 0.059483
$ws.Cells.Item(3,9).Value  = 0.1997488154364332
$ws.Cells.Item(3,10).Value = 0.1997488154364332
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.1790523333333333
$ws.Cells.Item(3,14).Value = 0.537157
$ws.Cells.Item(3,15).Value = 1
$ws.Cells.Item(3,16).Value = 1
$ws.Cells.Item(3,17).Value = 0.003550189981222222
$ws.Cells.Item(3,18).Value = 0.031951709831
$ws.Cells.Item(3,19).Value = 0.1997488154364332
$ws.Cells.Item(3,20).Value = 0.1997488154364332

# Remove the now-obsolete trailing rows 4 and 5 (FAPs->ECs, FAPs->Resolving-Mac)
# so the used range shrinks back down to A1:T3.
$ws.Range("A4:T5").EntireRow.Delete()
